$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31/32 coin identity swap (Filecoin <-> Toncoin) plus refreshed
# price/volume figures for every row, per the scraper run.

$ws.Range("D2").Value = '31.351.60'
$ws.Range("E2").Value = '  +3.21%  '
$ws.Range("D3").Value = '1.996.23'
$ws.Range("E3").Value = '  +6.71%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9993'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.8088'
$ws.Range("E5").Value = '  +71.74%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '254.68'
$ws.Range("E6").Value = '  +4.42%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9990'
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3526'
$ws.Range("E8").Value = '  +22.31%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.58'
$ws.Range("E9").Value = '  +16.43%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07024'
$ws.Range("E10").Value = '  +8.69%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8415'
$ws.Range("E11").Value = '  +16.23%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08124'
$ws.Range("E12").Value = '  +4.31%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '101.18'
$ws.Range("E13").Value = '  +5.24%  '
$ws.Range("D14").Value = '1.988.48'
$ws.Range("E14").Value = '  +6.61%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.540'
$ws.Range("E15").Value = '  +7.81%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '272.95'
$ws.Range("E16").Value = '  -3.12%  '
$ws.Range("D17").Value = '31.337.55'
$ws.Range("E17").Value = '  +3.19%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.00'
$ws.Range("E18").Value = '  +7.50%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007960'
$ws.Range("E19").Value = '  +6.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.810'
$ws.Range("E20").Value = '  +10.59%  '
$ws.Range("D21").Value = '2.251.98'
$ws.Range("E21").Value = '  +6.48%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9987'
$ws.Range("E22").Value = '  -0.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9991'
$ws.Range("E23").Value = '  -0.12%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.973'
$ws.Range("E24").Value = '  +11.74%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.826'
$ws.Range("E25").Value = '  +8.51%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1509'
$ws.Range("E26").Value = '  +56.98%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '164.10'
$ws.Range("E27").Value = '  +0.21%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.08'
$ws.Range("E28").Value = '  +7.16%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.231'
$ws.Range("E29").Value = '  +18.62%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.572'
$ws.Range("E30").Value = '  +5.86%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.357'
$ws.Range("E31").Value = '  +2.77%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.591'
$ws.Range("E32").Value = '  +8.53%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.342'
$ws.Range("E33").Value = '  +5.56%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05205'
$ws.Range("E34").Value = '  +8.09%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.216'
$ws.Range("E35").Value = '  +8.50%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7597'
$ws.Range("E36").Value = '  +10.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.774'
$ws.Range("E37").Value = '  +2.25%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02012'
$ws.Range("E38").Value = '  +6.45%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.914'
$ws.Range("E39").Value = '  +3.35%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.656'
$ws.Range("E40").Value = '  +6.86%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4763'
$ws.Range("E41").Value = '  +12.79%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '78.31'
$ws.Range("E42").Value = '  +4.76%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.105'
$ws.Range("E43").Value = '  +9.09%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8563'
$ws.Range("E44").Value = '  +3.49%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '104.42'
$ws.Range("E45").Value = '  +3.40%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9992'
$ws.Range("E46").Value = '  -0.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.958'
$ws.Range("E47").Value = '  +3.20%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.524'
$ws.Range("E48").Value = '  +8.12%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4376'
$ws.Range("E49").Value = '  +11.46%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.96'
$ws.Range("E50").Value = '  +4.72%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1204'
$ws.Range("E51").Value = '  +13.41%  '
